$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activités")

# Row 26: "22 mai" -> "22 mai 2020" (same row content otherwise)
$ws.Range("A26").Value = "22 mai 2020"

# Row 27 (was blank template row): 26 mai 2020, 08:00 - 12:00, Chat
$ws.Range("A27").Value = "26 mai 2020"
$ws.Range("B27").Value = 0.33333333333333331
$ws.Range("C27").Value = 0.5
$ws.Range("D27").Formula = "=IF(ISBLANK(C27), NOW(),C27)-IF(ISBLANK(B27),NOW(),B27)"
$ws.Range("E27").Value = "Chat"
$ws.Range("F27").Value = "Réalisation"
$ws.Range("G27").Value = "J'ai corrigé les bugs lors de l'envoi des utilisateurs connectés et non connectés et j'ai avancé dans la réalisation du chat."

# Row 28 (was blank template row): 26 mai 2020, 13:30 - 14:30, Chat
$ws.Range("A28").Value = "26 mai 2020"
$ws.Range("B28").Value = 0.5625
$ws.Range("C28").Value = 0.60416666666666663
$ws.Range("D28").Formula = "=IF(ISBLANK(C28), NOW(),C28)-IF(ISBLANK(B28),NOW(),B28)"
$ws.Range("E28").Value = "Chat"
$ws.Range("F28").Value = "Réalisation"
$ws.Range("G28").Value = "Les messages peuvent maintenant être envoyés, mais il reste des bugs d'affichage."

$ws.Rows(27).RowHeight = 45
$ws.Rows(28).RowHeight = 30

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G34").Select()
